$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 457
$ws.Range("I2").Value = 1190
$ws.Range("J2").Value = 4928
$ws.Range("K2").Value = 11
$ws.Range("L2").Value = 1331
$ws.Range("M2").Value = 78
$ws.Range("N2").Value = 882
$ws.Range("O2").Value = 4
$ws.Range("P2").Value = 14
$ws.Range("Q2").Value = 8
$ws.Range("R2").Value = 69
$ws.Range("S2").Value = 566
$ws.Range("T2").Value = 878
$ws.Range("U2").Value = 67
$ws.Range("V2").Value = 7755
$ws.Range("W2").Value = 4
$ws.Range("X2").Value = 7842
$ws.Range("Y2").Value = 11
$ws.Range("Z2").Value = 135
$ws.Range("AA2").Value = 47
